$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6900
$ws.Range("E2").Value = 246
$ws.Range("F2").Value = 246
$ws.Range("G2").Value = 65
$ws.Range("H2").Value = 43
$ws.Range("I2").Value = 43
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7207
$ws.Range("L2").Value = 5611
$ws.Range("M2").Value = 1597
$ws.Range("N2").Value = 1528
$ws.Range("O2").Value = 68
$ws.Range("P2").Value = 107
$ws.Range("Q2").Value = 546
$ws.Range("R2").Value = -331
$ws.Range("S2").Value = -204
$ws.Range("T2").Value = 272
$ws.Range("U2").Value = 274
$ws.Range("V2").Value = 3398
$ws.Range("W2").Value = 3.57
$ws.Range("X2").Value = 0.62
$ws.Range("Y2").Value = 2.84
$ws.Range("Z2").Value = 0.6
$ws.Range("AA2").Value = 351.43
$ws.Range("AB2").Value = 1744.77
$ws.Range("AC2").Value = 181
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 6415
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 0
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 13150240

# Row 3
$ws.Range("D3").Value = 7983
$ws.Range("E3").Value = 469
$ws.Range("F3").Value = 469
$ws.Range("G3").Value = 244
$ws.Range("H3").Value = 169
$ws.Range("I3").Value = 169
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 7482
$ws.Range("L3").Value = 5713
$ws.Range("M3").Value = 1769
$ws.Range("N3").Value = 1688
$ws.Range("O3").Value = 80
$ws.Range("P3").Value = 119
$ws.Range("Q3").Value = 649
$ws.Range("R3").Value = -406
$ws.Range("S3").Value = -180
$ws.Range("T3").Value = 323
$ws.Range("U3").Value = 326
$ws.Range("V3").Value = 3362
$ws.Range("W3").Value = 5.87
$ws.Range("X3").Value = 2.12
$ws.Range("Y3").Value = 10.49
$ws.Range("Z3").Value = 2.3
$ws.Range("AA3").Value = 323.04
$ws.Range("AB3").Value = 2029.94
$ws.Range("AC3").Value = 708
$ws.Range("AD3").ClearContents()
$ws.Range("AE3").Value = 7086
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").ClearContents()
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 15504900

# Row 4
$ws.Range("D4").Value = 7928
$ws.Range("E4").Value = 352
$ws.Range("F4").Value = 352
$ws.Range("G4").Value = 320
$ws.Range("H4").Value = 255
$ws.Range("I4").Value = 253
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 7963
$ws.Range("L4").Value = 5016
$ws.Range("M4").Value = 2947
$ws.Range("N4").Value = 2865
$ws.Range("O4").Value = 82
$ws.Range("P4").Value = 146
$ws.Range("Q4").Value = 354
$ws.Range("R4").Value = -496
$ws.Range("S4").Value = 90
$ws.Range("T4").Value = 481
$ws.Range("U4").Value = -127
$ws.Range("V4").Value = 2657
$ws.Range("W4").Value = 4.43
$ws.Range("X4").Value = 3.21
$ws.Range("Y4").Value = 11.11
$ws.Range("Z4").Value = 3.3
$ws.Range("AA4").Value = 170.22
$ws.Range("AB4").Value = 2191.48
$ws.Range("AC4").Value = 927
$ws.Range("AD4").Value = 19.36
$ws.Range("AE4").Value = 10803
$ws.Range("AF4").Value = 1.66
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 0.84
$ws.Range("AI4").Value = 15.73
$ws.Range("AJ4").Value = 29116822

# Row 5
$ws.Range("D5").Value = 7604
$ws.Range("E5").Value = 189
$ws.Range("F5").Value = 189
$ws.Range("G5").Value = 87
$ws.Range("H5").Value = 75
$ws.Range("I5").Value = 75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 8209
$ws.Range("L5").Value = 5313
$ws.Range("M5").Value = 2895
$ws.Range("N5").Value = 2882
$ws.Range("O5").Value = 13
$ws.Range("P5").Value = 146
$ws.Range("Q5").Value = 549
$ws.Range("R5").Value = -256
$ws.Range("S5").Value = -313
$ws.Range("T5").Value = 306
$ws.Range("U5").Value = 243
$ws.Range("V5").Value = 2881
$ws.Range("W5").Value = 2.49
$ws.Range("X5").Value = 0.98
$ws.Range("Y5").Value = 2.59
$ws.Range("Z5").Value = 0.92
$ws.Range("AA5").Value = 183.54
$ws.Range("AB5").Value = 2203.62
$ws.Range("AC5").Value = 256
$ws.Range("AD5").Value = 60.36
$ws.Range("AE5").Value = 10870
$ws.Range("AF5").Value = 1.42
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 0.97
$ws.Range("AI5").Value = 53.37
$ws.Range("AJ5").Value = 29116822

# Row 6
$ws.Range("D6").Value = 7254
$ws.Range("E6").Value = 230
$ws.Range("F6").Value = 230
$ws.Range("G6").Value = 42
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 24
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 7950
$ws.Range("L6").Value = 5074
$ws.Range("M6").Value = 2875
$ws.Range("N6").Value = 2863
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 146
$ws.Range("Q6").Value = 507
$ws.Range("R6").Value = -453
$ws.Range("S6").Value = -44
$ws.Range("T6").Value = 448
$ws.Range("U6").Value = 59
$ws.Range("V6").Value = 2977
$ws.Range("W6").Value = 3.17
$ws.Range("X6").Value = 0.34
$ws.Range("Y6").Value = 0.85
$ws.Range("Z6").Value = 0.3
$ws.Range("AA6").Value = 176.47
$ws.Range("AB6").Value = 2193.01
$ws.Range("AC6").Value = 84
$ws.Range("AD6").Value = 111.23
$ws.Range("AE6").Value = 10796
$ws.Range("AF6").Value = 0.86
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 1.61
$ws.Range("AI6").Value = 163.4
$ws.Range("AJ6").Value = 29116822

# Rows 7,8,9: remove all data columns D:AJ, keep A/B/C
$ws.Range("D7:AJ9").ClearContents()
